# "Responded to some comments": refresh the Scenario Analysis Timeframe
# model outputs for Net.Monetary.Benefit..Overall. (C), QALYs.Saved (D),
# Net.Monetary.Benefit..Healthcare. (E) and Net.Monetary.Benefit..Productivity. (F)
# across timeframes 1,2,3,4,5 (rows 2-6).
#
# The source values are stored as *text* (not numbers) in the workbook.
# Typing a numeric-looking value directly into Range.Value would make Excel
# coerce it to a real number (and/or tag the cell with a quote-prefix style
# if forced to text via NumberFormat), changing the cells type/style versus
# the original. To keep the cell a plain text (shared-string) cell exactly
# like the source - with no style change - write the new value as a text
# formula result, then flatten it to a static value with copy / paste-special
# values (xlPasteValues = -4163), mirroring how the numbers were authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Formula = '="912261137.409362"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

$ws.Range("C3").Formula = '="1145084169.00089"'
$ws.Range("C3").Copy()
$ws.Range("C3").PasteSpecial(-4163)

$ws.Range("C4").Formula = '="1092965034.95472"'
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("C5").Formula = '="947897913.909862"'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$ws.Range("C6").Formula = '="850894453.187894"'
$ws.Range("C6").Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range("D2").Formula = '="1092.67382125801"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="2424.53610411938"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("D4").Formula = '="3207.21197185642"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="3602.5245692532"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="3732.83737464552"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("E2").Formula = '="2975321.26930352"'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("E3").Formula = '="6322139.16841889"'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("E4").Formula = '="8273548.72198048"'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("E5").Formula = '="9260278.7945701"'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("E6").Formula = '="9586774.16410338"'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("F2").Formula = '="7949952.7462759"'
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

$ws.Range("F3").Formula = '="15350533.3239965"'
$ws.Range("F3").Copy()
$ws.Range("F3").PasteSpecial(-4163)

$ws.Range("F4").Formula = '="21547591.0756111"'
$ws.Range("F4").Copy()
$ws.Range("F4").PasteSpecial(-4163)

$ws.Range("F5").Formula = '="26762914.2908649"'
$ws.Range("F5").Copy()
$ws.Range("F5").PasteSpecial(-4163)

$ws.Range("F6").Formula = '="29500633.6890926"'
$ws.Range("F6").Copy()
$ws.Range("F6").PasteSpecial(-4163)

$excel.CutCopyMode = 0
